$d = $word.ActiveDocument

# 1. Update the hidden ID placeholder text in the first paragraph and drop
#    the trailing space-only run by replacing the whole "text + space" span
#    with just the new ID text.
$d.Content.Find.Execute("**ID__AFFARS_5301_topic_16__ID** ", $true, $false, $false, $false, $false, $true, 1, $false, "**ID__AFFARS_SUBPART_5301_6__ID**", 2) | Out-Null

# 2. Give that same (first) paragraph a paragraph border (top/left/bottom/right)
#    whose only attribute is a 5-twip space, and bump its left indent from
#    120 to 225 twips (6pt -> 11.25pt).
$p1 = $d.Paragraphs.First
$borders = $p1.Borders
$borders.DistanceFromTop = 5
$borders.DistanceFromLeft = 5
$borders.DistanceFromBottom = 5
$borders.DistanceFromRight = 5
$p1.LeftIndent = 11.25
